# Update TPM-derived values on Sheet1 (Edn3-Ednra LR-pairs) to reflect
# the new TPM-based recalculation referenced in the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.08840121110588733
$ws.Range("J2").Value = 0.08840121110588735
$ws.Range("M2").Value = 2.067959333333333
$ws.Range("N2").Value = 6.203878
$ws.Range("O2").Value = 0.03751906553627715
$ws.Range("P2").Value = 0.03751906553627715
$ws.Range("Q2").Value = 0.7362079983819999
$ws.Range("R2").Value = 6.625871985438001
$ws.Range("S2").Value = 0.003316730832968059
$ws.Range("T2").Value = 0.003316730832968059

# Row 3
$ws.Range("I3").Value = 0.08840121110588733
$ws.Range("J3").Value = 0.08840121110588735
$ws.Range("O3").Value = 0.1625861532004571
$ws.Range("P3").Value = 0.1625861532004571
$ws.Range("R3").Value = 28.712736373176
$ws.Range("S3").Value = 0.01437281285196775
$ws.Range("T3").Value = 0.01437281285196775

# Row 4
$ws.Range("I4").Value = 0.08840121110588733
$ws.Range("J4").Value = 0.08840121110588735
$ws.Range("O4").Value = 0.7998947812632657
$ws.Range("P4").Value = 0.7998947812632657
$ws.Range("S4").Value = 0.07071166742095153
$ws.Range("T4").Value = 0.07071166742095154

# Row 5
$ws.Range("M5").Value = 2.067959333333333
$ws.Range("N5").Value = 6.203878
$ws.Range("O5").Value = 0.03751906553627715
$ws.Range("P5").Value = 0.03751906553627715
$ws.Range("Q5").Value = 7.591822683235777
$ws.Range("R5").Value = 68.326404149122
$ws.Range("S5").Value = 0.03420233470330909
$ws.Range("T5").Value = 0.03420233470330909

# Row 6
$ws.Range("O6").Value = 0.1625861532004571
$ws.Range("P6").Value = 0.1625861532004571
$ws.Range("S6").Value = 0.1482133403484893
$ws.Range("T6").Value = 0.1482133403484893

# Row 7
$ws.Range("O7").Value = 0.7998947812632657
$ws.Range("P7").Value = 0.7998947812632657
$ws.Range("S7").Value = 0.7291831138423142
$ws.Range("T7").Value = 0.7291831138423143
